# Error Calculations and Plots
# Apply imputation/correction edits to the "missing_data" worksheet:
#  1. Remove the "RM 232" and "SC 92" rows (rows shift up by two).
#  2. Fill / clear a handful of column D (and a couple of column C) values
#     to match the corrected data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two rows that were dropped from the data set ---------
# Row 26 is "RM 232". After it is removed, the row that used to be 28
# ("SC 92") becomes row 27, so we delete row 27 next.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- 2. Column D corrections among rows 2-24 -----------------------------
$ws.Range("D2").Value2 = -13.5
$ws.Range("D6").ClearContents()
$ws.Range("D12").Value2 = -14.1
$ws.Range("D14").ClearContents()
$ws.Range("D20").Value2 = -14
$ws.Range("D21").Value2 = -14.3
$ws.Range("D23").ClearContents()
$ws.Range("D24").ClearContents()

# --- 3. Corrections among the shifted SC rows (now rows 26-33) ----------
$ws.Range("C26").Value2 = 10.8
$ws.Range("C27").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("C29").Value2 = 11.2
$ws.Range("C30").Value2 = 11.4
$ws.Range("C31").ClearContents()
$ws.Range("D31").Value2 = -13.7
$ws.Range("C32").ClearContents()
$ws.Range("D33").Value2 = -14.1

Write-Output "Applied missing-data corrections; used range is now $($ws.UsedRange.Address())"
